$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text, $styleSrcRow, $styleSrcCol) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $text
    $src = $ws.Cells.Item($styleSrcRow, $styleSrcCol)
    $src.Copy()
    $c.PasteSpecial(-4122)
}

function Set-NumCellWithStyle($row, $col, $val, $styleSrcRow, $styleSrcCol) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $val
    $src = $ws.Cells.Item($styleSrcRow, $styleSrcCol)
    $src.Copy()
    $c.PasteSpecial(-4122)
}

# --- Value updates ---
$ws.Cells.Item(14, 14).Value = -82.758620689655  # N14
Set-TextCell 15 7 "0" 14 3  # G15
Set-TextCell 15 8 "***.*" 14 5  # H15
$ws.Cells.Item(15, 13).Value = -22.222222222222  # M15
$ws.Cells.Item(15, 14).Value = -78.125  # N15
Set-TextCell 16 3 "0" 14 3  # C16
$ws.Cells.Item(16, 4).Value = 2  # D16
$ws.Cells.Item(16, 5).Value = -100  # E16
$ws.Cells.Item(16, 6).Value = 4  # F16
$ws.Cells.Item(16, 7).Value = 7  # G16
$ws.Cells.Item(16, 8).Value = -42.857142857142  # H16
$ws.Cells.Item(16, 9).Value = 94  # I16
$ws.Cells.Item(16, 10).Value = 132  # J16
$ws.Cells.Item(16, 11).Value = -28.787878787878  # K16
$ws.Cells.Item(16, 12).Value = -30.37037037037  # L16
$ws.Cells.Item(16, 13).Value = -47.486033519553  # M16
$ws.Cells.Item(16, 14).Value = -85.779122541603  # N16
$ws.Cells.Item(17, 3).Value = 3  # C17
$ws.Cells.Item(17, 4).Value = 4  # D17
$ws.Cells.Item(17, 5).Value = -25  # E17
$ws.Cells.Item(17, 6).Value = 15  # F17
$ws.Cells.Item(17, 7).Value = 20  # G17
$ws.Cells.Item(17, 8).Value = -25  # H17
$ws.Cells.Item(17, 9).Value = 191  # I17
$ws.Cells.Item(17, 10).Value = 233  # J17
$ws.Cells.Item(17, 11).Value = -18.025751072961  # K17
$ws.Cells.Item(17, 12).Value = -3.045685279187  # L17
$ws.Cells.Item(17, 13).Value = 31.724137931034  # M17
$ws.Cells.Item(17, 14).Value = -66.608391608391  # N17
$ws.Cells.Item(18, 3).Value = 4  # C18
$ws.Cells.Item(18, 4).Value = 2  # D18
$ws.Cells.Item(18, 5).Value = 100  # E18
$ws.Cells.Item(18, 6).Value = 9  # F18
$ws.Cells.Item(18, 7).Value = 9  # G18
$ws.Cells.Item(18, 8).Value = 0  # H18
$ws.Cells.Item(18, 9).Value = 56  # I18
$ws.Cells.Item(18, 10).Value = 71  # J18
$ws.Cells.Item(18, 11).Value = -21.12676056338  # K18
$ws.Cells.Item(18, 12).Value = -35.632183908046  # L18
$ws.Cells.Item(18, 13).Value = -36.363636363636  # M18
$ws.Cells.Item(18, 14).Value = -91.591591591591  # N18
$ws.Cells.Item(19, 3).Value = 5  # C19
$ws.Cells.Item(19, 4).Value = 7  # D19
$ws.Cells.Item(19, 5).Value = -28.571428571428  # E19
$ws.Cells.Item(19, 6).Value = 39  # F19
$ws.Cells.Item(19, 7).Value = 38  # G19
$ws.Cells.Item(19, 8).Value = 2.631578947368  # H19
$ws.Cells.Item(19, 9).Value = 309  # I19
$ws.Cells.Item(19, 10).Value = 393  # J19
$ws.Cells.Item(19, 11).Value = -21.374045801526  # K19
$ws.Cells.Item(19, 12).Value = 5.821917808219  # L19
$ws.Cells.Item(19, 13).Value = 32.618025751073  # M19
$ws.Cells.Item(19, 14).Value = -14.404432132964  # N19
Set-TextCell 20 6 "0" 14 3  # F20
$ws.Cells.Item(20, 8).Value = -100  # H20
$ws.Cells.Item(20, 10).Value = 47  # J20
$ws.Cells.Item(20, 11).Value = -53.191489361702  # K20
$ws.Cells.Item(20, 12).Value = -59.259259259259  # L20
$ws.Cells.Item(20, 13).Value = 29.411764705882  # M20
$ws.Cells.Item(20, 14).Value = -80.18018018018  # N20
$ws.Cells.Item(21, 3).Value = 12  # C21
$ws.Cells.Item(21, 4).Value = 15  # D21
$ws.Cells.Item(21, 5).Value = -20  # E21
$ws.Cells.Item(21, 7).Value = 75  # G21
$ws.Cells.Item(21, 8).Value = -10.666666666666  # H21
$ws.Cells.Item(21, 9).Value = 684  # I21
$ws.Cells.Item(21, 10).Value = 888  # J21
$ws.Cells.Item(21, 11).Value = -22.972972972973  # K21
$ws.Cells.Item(21, 12).Value = -11.627906976744  # L21
$ws.Cells.Item(21, 13).Value = 1.183431952662  # M21
$ws.Cells.Item(21, 14).Value = -71.875  # N21
Set-TextCell 22 3 "0" 14 3  # C22
Set-TextCell 22 4 "0" 14 3  # D22
Set-TextCell 22 5 "***.*" 14 5  # E22
$ws.Cells.Item(22, 6).Value = 2  # F22
$ws.Cells.Item(22, 8).Value = 0  # H22
$ws.Cells.Item(22, 9).Value = 14  # I22
$ws.Cells.Item(22, 11).Value = -12.5  # K22
$ws.Cells.Item(22, 12).Value = -30  # L22
$ws.Cells.Item(22, 13).Value = 75  # M22
$ws.Cells.Item(23, 6).Value = 2  # F23
$ws.Cells.Item(23, 7).Value = 1  # G23
$ws.Cells.Item(23, 8).Value = 100  # H23
$ws.Cells.Item(23, 9).Value = 46  # I23
$ws.Cells.Item(23, 11).Value = -29.230769230769  # K23
$ws.Cells.Item(23, 12).Value = -32.35294117647  # L23
$ws.Cells.Item(23, 13).Value = 31.428571428571  # M23
$ws.Cells.Item(24, 3).Value = 29  # C24
$ws.Cells.Item(24, 4).Value = 25  # D24
$ws.Cells.Item(24, 5).Value = 16  # E24
$ws.Cells.Item(24, 6).Value = 111  # F24
$ws.Cells.Item(24, 7).Value = 108  # G24
$ws.Cells.Item(24, 8).Value = 2.777777777777  # H24
$ws.Cells.Item(24, 9).Value = 1168  # I24
$ws.Cells.Item(24, 10).Value = 964  # J24
$ws.Cells.Item(24, 11).Value = 21.161825726141  # K24
$ws.Cells.Item(24, 12).Value = 19.794871794871  # L24
$ws.Cells.Item(24, 13).Value = 57.200538358008  # M24
$ws.Cells.Item(25, 3).Value = 16  # C25
$ws.Cells.Item(25, 4).Value = 16  # D25
$ws.Cells.Item(25, 5).Value = 0  # E25
$ws.Cells.Item(25, 6).Value = 62  # F25
$ws.Cells.Item(25, 7).Value = 58  # G25
$ws.Cells.Item(25, 8).Value = 6.896551724137  # H25
$ws.Cells.Item(25, 9).Value = 747  # I25
$ws.Cells.Item(25, 10).Value = 588  # J25
$ws.Cells.Item(25, 11).Value = 27.04081632653  # K25
$ws.Cells.Item(25, 12).Value = 29.91304347826  # L25
$ws.Cells.Item(26, 3).Value = 5  # C26
$ws.Cells.Item(26, 4).Value = 9  # D26
$ws.Cells.Item(26, 5).Value = -44.444444444444  # E26
$ws.Cells.Item(26, 6).Value = 29  # F26
$ws.Cells.Item(26, 7).Value = 39  # G26
$ws.Cells.Item(26, 8).Value = -25.641025641025  # H26
$ws.Cells.Item(26, 9).Value = 332  # I26
$ws.Cells.Item(26, 10).Value = 340  # J26
$ws.Cells.Item(26, 11).Value = -2.35294117647  # K26
$ws.Cells.Item(26, 12).Value = 6.070287539936  # L26
$ws.Cells.Item(26, 13).Value = -15.089514066496  # M26
Set-TextCell 27 7 "0" 14 3  # G27
Set-TextCell 27 8 "***.*" 14 5  # H27
Set-TextCell 28 3 "0" 14 3  # C28
Set-TextCell 28 4 "0" 14 3  # D28
Set-TextCell 28 5 "***.*" 14 5  # E28
$ws.Cells.Item(28, 6).Value = 4  # F28
$ws.Cells.Item(28, 7).Value = 2  # G28
$ws.Cells.Item(28, 8).Value = 100  # H28
$ws.Cells.Item(28, 9).Value = 31  # I28
$ws.Cells.Item(28, 11).Value = -11.428571428571  # K28
$ws.Cells.Item(28, 12).Value = 10.714285714285  # L28
Set-TextCell 29 7 "0" 14 3  # G29
Set-TextCell 29 8 "***.*" 14 5  # H29
$ws.Cells.Item(29, 14).Value = -90.588235294117  # N29
Set-TextCell 30 7 "0" 14 3  # G30
Set-TextCell 30 8 "***.*" 14 5  # H30
$ws.Cells.Item(30, 14).Value = -90.90909090909  # N30
Set-NumCellWithStyle 31 3 1 14 9  # C31 (text->num)
Set-NumCellWithStyle 31 6 1 14 9  # F31 (text->num)
$ws.Cells.Item(31, 9).Value = 3  # I31
$ws.Cells.Item(31, 11).Value = 50  # K31
$ws.Cells.Item(31, 12).Value = -25  # L31
